# Update "想去人数" (column F) values on sheets "展览" and "全部类型"
# to reflect refreshed counts as generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 1025
$ws1.Range("F8").Value  = 2099
$ws1.Range("F12").Value = 1666
$ws1.Range("F13").Value = 392
$ws1.Range("F19").Value = 624
$ws1.Range("F20").Value = 710
$ws1.Range("F21").Value = 597
$ws1.Range("F22").Value = 12202
$ws1.Range("F23").Value = 12240
$ws1.Range("F25").Value = 697
$ws1.Range("F27").Value = 29
$ws1.Range("F29").Value = 355
$ws1.Range("F30").Value = 1918
$ws1.Range("F31").Value = 192

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 1025
$ws4.Range("F9").Value  = 2099
$ws4.Range("F10").Value = 1099
$ws4.Range("F13").Value = 1666
$ws4.Range("F14").Value = 392
$ws4.Range("F23").Value = 624
$ws4.Range("F24").Value = 710
$ws4.Range("F25").Value = 597
$ws4.Range("F26").Value = 12202
$ws4.Range("F27").Value = 12240
$ws4.Range("F29").Value = 697
$ws4.Range("F31").Value = 29
$ws4.Range("F33").Value = 355
$ws4.Range("F34").Value = 1918
$ws4.Range("F37").Value = 192
